$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "emanuel"
$ws.Range("B6").Value = "ema"
$ws.Range("C6").Value = "Cliente"

$ws.Range("A7").Value = "igna"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "123"
$ws.Range("C7").Value = "Cliente"
